$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 2 (horse-racing AI gig) is replaced in place by the freshly
# scraped record, and the old row 3 (Go engineer gig) is dropped entirely
# since only one new listing was appended this run.
$ws.Range("A2").Value = "2025-12-14 06:27:09"
$ws.Range("B2").Value = "注目 限定公開 PR 限定公開の仕事"
$ws.Range("D2").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5450323"
$ws.Range("G2").Value = 13
$ws.Range("H2").ClearContents()

# Drop the old second record (row 3) entirely.
$ws.Rows.Item(3).Delete()

# Rebuild the hyperlink list: clear all existing links (engine quirk -
# deleting any one cell's Hyperlinks collection clears the whole sheet)
# then re-add only the link for the remaining F2 cell with the new URL.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5450323")
# Hyperlinks.Add silently tacks on a fresh (duplicate) "Hyperlink" cell
# style instead of reusing the workbook's existing one - reapply the
# named style explicitly so F2 keeps pointing at the original style.
$ws.Range("F2").Style = "Hyperlink"

# Column widths were retuned for the new (shorter) title/price text.
# ColumnWidth is in character units; the engine stores width = input +
# 5/6, so subtract 5/6 here to land exactly on the target stored widths
# of 20 and 26.
$ws.Columns.Item(2).ColumnWidth = 19.166666666666668
$ws.Columns.Item(4).ColumnWidth = 25.166666666666668
